$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
Write-Host "Col A width (orig):" $ws.Columns.Item(1).ColumnWidth
$ws.Columns.Item(1).ColumnWidth = 24.92
Write-Host "Col A width (set 24.92):" $ws.Columns.Item(1).ColumnWidth
